$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# machineindex -> machineno: replace numeric machine index codes with
# their string machine-number equivalents in column J.
$ws.Range("J2").Value = "O2310"
$ws.Range("J4").Value = "C2010"
